$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 3472503.2
$ws.Range("I53").Value = 5555797
$ws.Range("J53").Value = 346.91666
$ws.Range("K53").Value = 5555797
$ws.Range("L53").Value = 346.91666
$ws.Range("M53").Value = -5555160
$ws.Range("N53").Value = -1620.91666

$ws.Range("H87").Value = 11560.974
$ws.Range("J87").Value = 11560.974
$ws.Range("L87").Value = 11560.974
$ws.Range("N87").Value = -14056.974

$ws.Range("H90").Value = 11560.974
$ws.Range("J90").Value = 11560.974
$ws.Range("L90").Value = 34682.922
$ws.Range("N90").Value = -47162.922

$ws.Range("H107").Value = 723.0417
$ws.Range("I107").Value = 338.41177
$ws.Range("K107").Value = 338.41177
$ws.Range("M107").Value = 1581.58823

$ws.Range("H111").Value = 566.6667
$ws.Range("I111").Value = 540
$ws.Range("J111").Value = 700
$ws.Range("K111").Value = 1620
$ws.Range("L111").Value = 2100
$ws.Range("M111").Value = 1447
$ws.Range("N111").Value = -8234

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 740.1177
$ws.Range("I74").Value = 628.5769
$ws.Range("J74").Value = 1102.625
$ws.Range("K74").Value = 628.5769
$ws.Range("L74").Value = 1102.625
$ws.Range("M74").Value = 245.4231
$ws.Range("N74").Value = -2850.625

$ws.Range("H77").Value = 740.1177
$ws.Range("I77").Value = 628.5769
$ws.Range("J77").Value = 1102.625
$ws.Range("K77").Value = 3142.8845
$ws.Range("L77").Value = 5513.125
$ws.Range("M77").Value = 1225.1155
$ws.Range("N77").Value = -14249.125

$ws.Range("H122").Value = 7838273
$ws.Range("I122").Value = 34203.348
$ws.Range("J122").Value = 27782006
$ws.Range("K122").Value = 102610.044
$ws.Range("L122").Value = 83346018
$ws.Range("M122").Value = -100160.044
$ws.Range("N122").Value = -83350918

$ws.Range("H132").Value = 29471706
$ws.Range("I132").Value = 35715180
$ws.Range("J132").Value = 335501.34
$ws.Range("K132").Value = 107145540
$ws.Range("L132").Value = 1006504.02
$ws.Range("M132").Value = -107143010
$ws.Range("N132").Value = -1011564.02

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2021.1154
$ws.Range("I99").Value = 1149.2
$ws.Range("K99").Value = 1149.2
$ws.Range("M99").Value = 348.8

$ws.Range("H107").Value = 3164.8
$ws.Range("I107").Value = 3270.3333
$ws.Range("J107").Value = 3006.5
$ws.Range("K107").Value = 3270.3333
$ws.Range("L107").Value = 3006.5
$ws.Range("M107").Value = -1350.3333
$ws.Range("N107").Value = -6846.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1051732.6
$ws.Range("I31").Value = 1154392.8
$ws.Range("J31").Value = 4600
$ws.Range("K31").Value = 1154392.8
$ws.Range("L31").Value = 4600
$ws.Range("M31").Value = -1154097.8
$ws.Range("N31").Value = -5190

$ws.Range("H34").Value = 1051732.6
$ws.Range("I34").Value = 1154392.8
$ws.Range("J34").Value = 4600
$ws.Range("K34").Value = 1154392.8
$ws.Range("L34").Value = 4600
$ws.Range("M34").Value = -1154190.8
$ws.Range("N34").Value = -5004

$ws.Range("H58").Value = 1085
$ws.Range("I58").Value = 956.1539
$ws.Range("J58").Value = 1237.2727
$ws.Range("K58").Value = 956.1539
$ws.Range("L58").Value = 1237.2727
$ws.Range("M58").Value = -753.1539
$ws.Range("N58").Value = -1643.2727

$ws.Range("H81").Value = 29780
$ws.Range("J81").Value = 29780
$ws.Range("L81").Value = 29780
$ws.Range("N81").Value = -31776

$ws.Range("H84").Value = 29780
$ws.Range("J84").Value = 29780
$ws.Range("L84").Value = 89340
$ws.Range("N84").Value = -99324

$ws.Range("H88").Value = 18000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 18000
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 18000
$ws.Range("N88").Value = -18812

$ws.Range("H91").Value = 18000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 18000
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 18000
$ws.Range("N91").Value = -20808

$ws.Range("H132").Value = 4547750
$ws.Range("I132").Value = 1962.081
$ws.Range("J132").Value = 28575486
$ws.Range("K132").Value = 5886.242999999999
$ws.Range("L132").Value = 85726458
$ws.Range("M132").Value = -3356.242999999999
$ws.Range("N132").Value = -85731518

$ws.Range("H134").Value = 13595.75
$ws.Range("I134").Value = 6942.8945
$ws.Range("J134").Value = 140000
$ws.Range("K134").Value = 20828.6835
$ws.Range("L134").Value = 420000
$ws.Range("M134").Value = -18293.6835
$ws.Range("N134").Value = -425070

$ws.Range("H136").Value = 1085
$ws.Range("I136").Value = 956.1539
$ws.Range("J136").Value = 1237.2727
$ws.Range("K136").Value = 2868.4617
$ws.Range("L136").Value = 3711.8181
$ws.Range("M136").Value = -318.4616999999998
$ws.Range("N136").Value = -8811.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 73.111115
$ws.Range("I8").Value = 73.111115
$ws.Range("K8").Value = 219.333345
$ws.Range("M8").Value = -80.33334500000001

$ws.Range("H14").Value = 117
$ws.Range("I14").Value = 117
$ws.Range("K14").Value = 351
$ws.Range("M14").Value = -178

$ws.Range("H75").Value = 841.4286
$ws.Range("J75").Value = 841.4286
$ws.Range("L75").Value = 2524.2858
$ws.Range("N75").Value = -4520.2858

$ws.Range("H78").Value = 841.4286
$ws.Range("J78").Value = 841.4286
$ws.Range("L78").Value = 7572.8574
$ws.Range("N78").Value = -17556.8574

$ws.Range("H80").Value = 2411.3333
$ws.Range("I80").Value = 1402
$ws.Range("J80").Value = 2537.5
$ws.Range("K80").Value = 4206
$ws.Range("L80").Value = 7612.5
$ws.Range("M80").Value = -3270
$ws.Range("N80").Value = -9484.5

$ws.Range("H83").Value = 2411.3333
$ws.Range("I83").Value = 1402
$ws.Range("J83").Value = 2537.5
$ws.Range("K83").Value = 12618
$ws.Range("L83").Value = 22837.5
$ws.Range("M83").Value = -7938
$ws.Range("N83").Value = -32197.5

$ws.Range("H110").Value = 62505624
$ws.Range("I110").Value = 100002984
$ws.Range("K110").Value = 300008952
$ws.Range("M110").Value = -300004862

$ws.Range("H114").Value = 111111820
$ws.Range("I114").Value = 125000430
$ws.Range("J114").Value = 3000
$ws.Range("K114").Value = 375001290
$ws.Range("L114").Value = 9000
$ws.Range("M114").Value = -374998036
$ws.Range("N114").Value = -15508

$ws.Range("H120").Value = 13216.667
$ws.Range("J120").Value = 19000
$ws.Range("L120").Value = 57000
$ws.Range("N120").Value = -66676

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 7005
$ws.Range("I18").Value = 7005
$ws.Range("K18").Value = 7005
$ws.Range("M18").Value = -6712

$ws.Range("H102").Value = 1309.5883
$ws.Range("I102").Value = 1375.4546
$ws.Range("J102").Value = 1188.8334
$ws.Range("K102").Value = 1375.4546
$ws.Range("L102").Value = 1188.8334
$ws.Range("M102").Value = 246.5454
$ws.Range("N102").Value = -4432.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 22500
$ws.Range("I25").Value = 22500
$ws.Range("K25").Value = 22500
$ws.Range("M25").Value = -22270

$ws.Range("H136").Value = 50053140
$ws.Range("I136").Value = 91875.27
$ws.Range("J136").Value = 111116910
$ws.Range("K136").Value = 275625.81
$ws.Range("L136").Value = 333350730
$ws.Range("M136").Value = -273075.81
$ws.Range("N136").Value = -333355830

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1222.1666
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 2466.5
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 7399.5
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -11739.5

$ws.Range("H136").Value = 47974.953
$ws.Range("I136").Value = 59164.35
$ws.Range("J136").Value = 420
$ws.Range("K136").Value = 177493.05
$ws.Range("L136").Value = 1260
$ws.Range("M136").Value = -174943.05
$ws.Range("N136").Value = -6360
